$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
# A8: "Volume 32   Number  37" -> "...38"
$ws.Range("A8").Characters(21,2).Text = "38"

# C9: "Report Covering the Week  9/8/2025  Through  9/14/2025"
#  -> "Report Covering the Week  9/15/2025  Through  9/21/2025"
# Replace the later date first so the earlier substring offset stays valid.
$ws.Range("C9").Characters(46,9).Text = "9/21/2025"
$ws.Range("C9").Characters(27,8).Text = "9/15/2025"

# --- Numeric cell updates (weekly crime-stat figures) ---
$numericUpdates = @{
    "D15" = 1
    "F15" = 3
    "G15" = 4
    "H15" = -25
    "J15" = 18
    "K15" = 50
    "L15" = 125
    "N15" = -15.625
    "F16" = 6
    "G16" = 16
    "H16" = -62.5
    "I16" = 70
    "J16" = 99
    "K16" = -29.292929292929
    "L16" = -7.894736842105
    "M16" = -19.540229885057
    "N16" = -85.324947589098
    "C17" = 10
    "D17" = 7
    "E17" = 42.857142857142
    "F17" = 30
    "H17" = -18.918918918918
    "I17" = 277
    "J17" = 260
    "K17" = 6.538461538461
    "L17" = 42.783505154639
    "M17" = 151.818181818182
    "N17" = -27.486910994764
    "C18" = 2
    "E18" = 100
    "F18" = 7
    "G18" = 4
    "H18" = 75
    "I18" = 51
    "J18" = 52
    "K18" = -1.923076923076
    "L18" = -3.77358490566
    "M18" = -49.504950495049
    "N18" = -91.399662731871
    "C19" = 4
    "D19" = 1
    "E19" = 300
    "G19" = 12
    "H19" = 33.333333333333
    "I19" = 125
    "J19" = 122
    "K19" = 2.459016393442
    "L19" = -10.714285714285
    "M19" = 76.056338028169
    "N19" = -45.887445887445
    "D20" = 3
    "E20" = -100
    "F20" = 2
    "G20" = 8
    "H20" = -75
    "I20" = 26
    "J20" = 55
    "K20" = -52.727272727272
    "L20" = -46.938775510204
    "M20" = -59.375
    "N20" = -92.307692307692
    "C21" = 17
    "D21" = 15
    "E21" = 13.333333333333
    "F21" = 64
    "G21" = 81
    "H21" = -20.987654320987
    "I21" = 579
    "J21" = 608
    "K21" = -4.769736842105
    "L21" = 10.496183206106
    "M21" = 28.953229398663
    "N21" = -71.974830590513
    "D23" = 1
    "E23" = -100
    "F23" = 6
    "H23" = -25
    "I23" = 58
    "J23" = 63
    "K23" = -7.936507936507
    "L23" = -9.375
    "M23" = 123.076923076923
    "C24" = 14
    "D24" = 26
    "E24" = -46.153846153846
    "G24" = 67
    "H24" = 7.462686567164
    "I24" = 607
    "J24" = 476
    "K24" = 27.521008403361
    "L24" = 43.498817966903
    "M24" = 122.344322344322
    "D25" = 2
    "E25" = 50
    "F25" = 11
    "G25" = 3
    "H25" = 266.666666666667
    "I25" = 103
    "J25" = 36
    "K25" = 186.111111111111
    "L25" = 43.055555555555
    "C26" = 4
    "D26" = 12
    "E26" = -66.666666666666
    "F26" = 38
    "G26" = 32
    "H26" = 18.75
    "I26" = 353
    "J26" = 311
    "K26" = 13.504823151125
    "L26" = 12.063492063492
    "M26" = 17.666666666666
    "D27" = 2
    "F27" = 3
    "G27" = 6
    "H27" = -50
    "J27" = 29
    "K27" = 10.344827586206
    "L27" = 33.333333333333
    "C28" = 1
    "I28" = 29
    "K28" = -9.375
    "L28" = -12.121212121212
    "G29" = 2
    "H29" = -100
    "J29" = 13
    "K29" = -7.692307692307
    "M29" = -33.333333333333
    "G30" = 2
    "H30" = -100
    "J30" = 11
    "K30" = -9.090909090909
    "M30" = -41.176470588235
}
foreach ($ref in $numericUpdates.Keys) {
    $ws.Range($ref).Value = $numericUpdates[$ref]
}

# --- Cells that must become literal text placeholders ("0" / "***.*") ---
# A leading apostrophe forces Excel to store these as text instead of numbers,
# matching the "0" / "***.*" placeholder strings used elsewhere in the sheet.
$textUpdates = @{
    "C20" = "0"
    "C22" = "0"
    "C23" = "0"
    "D28" = "0"
    "E28" = "***.*"
    "F29" = "0"
    "F30" = "0"
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = "'" + $textUpdates[$ref]
}
